$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s2 = '$ 27.333 CLP 19-10-20'
$s3 = '$ 27.339 CLP 20-10-20'
$s4 = '$ 27.344 CLP 21-10-20'
$s5 = '$ 27.354 CLP 23-10-20'

$ws.Range("A7").Value = $s2
$ws.Range("A8").Value = $s2
$ws.Range("A9").Value = $s2
$ws.Range("A10").Value = $s2
$ws.Range("A11").Value = $s2
$ws.Range("A12").Value = $s2

$ws.Range("A13").Value = $s3
$ws.Range("A14").Value = $s3

$ws.Range("A15").Value = $s4
$ws.Range("A16").Value = $s4

$ws.Range("A17").Value = $s5
$ws.Range("A18").Value = $s5
$ws.Range("A19").Value = $s5
